# Re-label/normalize the macbook model names on the "prices_usados" sheet
# and fix two "1 AOA" placeholder prices, matching the author's
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prices_usados")
$ws.Activate()

$ws.Range("A47").Value = "Macbook Air 2015"
$ws.Range("A49").Value = "Macbook pro 2015"
$ws.Range("A51").Value = "Macbook Pro  Touch Bar 2018"
$ws.Range("A52").Value = "macbook pro touch bar 2018"
$ws.Range("A53").Value = "Macbook Pro  Touch Bar2018"
$ws.Range("A54").Value = "Macbook Pro  Touch Bar 2019"
$ws.Range("A55").Value = "macbook pro  touch bar 2019"
$ws.Range("A56").Value = "macbook pro  touch bar 2019"
$ws.Range("A57").Value = "macbook pro  touch bar 2019"
$ws.Range("A58").Value = "Macbook Pro  Touch Bar2020"
$ws.Range("A59").Value = "Macbook Pro  Touch Bar 2020"
$ws.Range("A60").Value = "Macbook Pro  Touch Bar 2020"
$ws.Range("A61").Value = "Macbook Pro  Touch Bar 2020"
$ws.Range("A62").Value = "iPad 10"
$ws.Range("A63").Value = "iPad 10"
$ws.Range("A64").Value = "iPad 11"
$ws.Range("A65").Value = "iPad 11"
$ws.Range("A66").Value = "iPad Pro"
$ws.Range("A67").Value = "iPad Air"

$ws.Range("E66").Value = 899000
$ws.Range("E67").Value = 755000

$ws.Range("A47").Select()
